$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.126.90'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.526.97'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.45%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.85'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.08'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.527.02'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.493'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.98%  '

$ws.Range("E10").Value = '  +2.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.82'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.91%  '

$ws.Range("E12").Value = '  +2.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.122.86'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000185'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.09'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.525.49'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.55%  '

$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.217.70'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.29'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.96'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '392.08'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.571'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.668.46'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.80'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000115'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +6.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +6.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.29'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.17'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.539.06'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.145'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("E36").Value = '  +6.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.94'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.55'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.56%  '

$ws.Range("E39").Value = '  +3.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.97'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0798'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.821'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.96'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +14.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.73'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.24%  '

$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.41'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.67'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.63%  '

$ws.Range("E48").Value = '  +4.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.79'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.374.05'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '301.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.00%  '
